$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new current-week row at row 16, pushing the
# previously-existing rows 16 and 17 down to rows 17 and 18 (history is
# preserved, a fresh observation is recorded up top of the recent data).
$ws.Rows(16).Insert()

$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 44644
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100112003
$ws.Cells.Item(16, 7).Value = "Ajo"
$ws.Cells.Item(16, 8).Value = "Chino"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 20000
$ws.Cells.Item(16, 12).Value = 21000
$ws.Cells.Item(16, 13).Value = 20500
$ws.Cells.Item(16, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(16, 15).Value = "China"
$ws.Cells.Item(16, 16).Value = 2050
$ws.Cells.Item(16, 17).Value = 10
$ws.Cells.Item(16, 18).Value = "Hortaliza"
